$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.004.79"
$ws.Range("E2").Value = "  +1.60%  "
$ws.Range("D3").Value = "2.498.91"
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").Value = "'534.22"
$ws.Range("E5").Value = "  +1.40%  "
$ws.Range("D6").Value = "'136.47"
$ws.Range("E6").Value = "  +2.53%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("E8").Value = "  +2.44%  "
$ws.Range("D9").Value = "2.515.64"
$ws.Range("E9").Value = "  +2.69%  "
$ws.Range("E10").Value = "  +3.41%  "
$ws.Range("E11").Value = "  -1.62%  "
$ws.Range("D12").Value = "'5.41"
$ws.Range("E12").Value = "  +2.09%  "
$ws.Range("E13").Value = "  +2.42%  "
$ws.Range("D14").Value = "2.964.85"
$ws.Range("E14").Value = "  +2.90%  "
$ws.Range("D15").Value = "'22.95"
$ws.Range("E15").Value = "  +2.30%  "
$ws.Range("D16").Value = "58.901.64"
$ws.Range("E16").Value = "  +1.57%  "
$ws.Range("E17").Value = "  +0.61%  "
$ws.Range("D18").Value = "2.515.55"
$ws.Range("E18").Value = "  +2.59%  "
$ws.Range("D19").Value = "'11.08"
$ws.Range("E19").Value = "  +4.38%  "
$ws.Range("E20").Value = "  +2.50%  "
$ws.Range("D21").Value = "'322.96"
$ws.Range("E21").Value = "  +1.35%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "'5.99"
$ws.Range("E22").Value = "  +5.50%  "
$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  +0.41%  "
$ws.Range("D24").Value = "'65.28"
$ws.Range("E24").Value = "  +5.07%  "
$ws.Range("E25").Value = "  +3.98%  "
$ws.Range("E26").Value = "  +0.46%  "
$ws.Range("D27").Value = "'0.997"
$ws.Range("E27").Value = "  +1.61%  "
$ws.Range("E28").Value = "  +1.42%  "
$ws.Range("D29").Value = "0.0₃0767"
$ws.Range("E29").Value = "  +3.10%  "
$ws.Range("D30").Value = "'6.59"
$ws.Range("E30").Value = "  +1.71%  "
$ws.Range("D31").Value = "'171.98"
$ws.Range("E31").Value = "  +5.73%  "
$ws.Range("E32").Value = "  +0.44%  "
$ws.Range("E33").Value = "  +10.04%  "
$ws.Range("D34").Value = "'0.998"
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("D35").Value = "'18.35"
$ws.Range("E35").Value = "  +1.43%  "
$ws.Range("E36").Value = "  +1.09%  "
$ws.Range("D37").Value = "'4.05"
$ws.Range("E37").Value = "  +1.63%  "
$ws.Range("D38").Value = "'1.53"
$ws.Range("E38").Value = "  +0.75%  "
$ws.Range("D39").Value = "'36.89"
$ws.Range("E39").Value = "  +1.71%  "
$ws.Range("E40").Value = "  +3.51%  "
$ws.Range("E41").Value = "  +2.21%  "
$ws.Range("D42").Value = "'283.52"
$ws.Range("E42").Value = "  +4.35%  "
$ws.Range("E43").Value = "  +1.45%  "
$ws.Range("D44").Value = "'0.997"
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("E45").Value = "  +4.08%  "
$ws.Range("D46").Value = "'130.56"
$ws.Range("E46").Value = "  +9.06%  "
$ws.Range("D47").Value = "'10.89"
$ws.Range("E47").Value = "  +0.56%  "
$ws.Range("D48").Value = "'0.0921"
$ws.Range("E48").Value = "  +0.45%  "
$ws.Range("E49").Value = "  +0.37%  "
$ws.Range("E50").Value = "  +1.44%  "
$ws.Range("D51").Value = "'17.31"
$ws.Range("E51").Value = "  +2.98%  "
